$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (interested count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5916
$ws1.Range("F3").Value = 560
$ws1.Range("F4").Value = 1105
$ws1.Range("F5").Value = 1058
$ws1.Range("F6").Value = 855
$ws1.Range("F7").Value = 89
$ws1.Range("F9").Value = 618
$ws1.Range("F10").Value = 63
$ws1.Range("F11").Value = 30
$ws1.Range("F13").Value = 2064
$ws1.Range("F14").Value = 1525
$ws1.Range("F15").Value = 1132
$ws1.Range("F17").Value = 211
$ws1.Range("F18").Value = 436
$ws1.Range("F19").Value = 661
$ws1.Range("F24").Value = 3740
$ws1.Range("F25").Value = 201
$ws1.Range("F30").Value = 531
$ws1.Range("F33").Value = 24
$ws1.Range("F35").Value = 331
$ws1.Range("F36").Value = 856
$ws1.Range("F37").Value = 109
$ws1.Range("F39").Value = 88
$ws1.Range("F40").Value = 92
# Row 16: ticket became available again -> F count updates and G now shows a numeric
# min price (218) instead of the "已售罄" (sold out) text
$ws1.Range("F16").Value = 307
$ws1.Range("G16").Value = 218

# Sheet "演出" (Performance) - update "想去人数" column F
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 756

# Sheet "全部类型" (All Types) - update "想去人数" column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5916
$ws4.Range("F4").Value = 560
$ws4.Range("F5").Value = 1105
$ws4.Range("F7").Value = 756
$ws4.Range("F8").Value = 1058
$ws4.Range("F9").Value = 855
$ws4.Range("F12").Value = 89
$ws4.Range("F14").Value = 618
$ws4.Range("F15").Value = 63
$ws4.Range("F16").Value = 30
$ws4.Range("F19").Value = 2065
$ws4.Range("F20").Value = 1525
$ws4.Range("F21").Value = 1132
$ws4.Range("F23").Value = 211
$ws4.Range("F24").Value = 436
$ws4.Range("F26").Value = 661
$ws4.Range("F30").Value = 3740
$ws4.Range("F31").Value = 201
$ws4.Range("F36").Value = 531
$ws4.Range("F39").Value = 24
$ws4.Range("F41").Value = 331
$ws4.Range("F42").Value = 856
$ws4.Range("F43").Value = 109
$ws4.Range("F45").Value = 88
$ws4.Range("F46").Value = 92
# Row 22: same event as row 16 on "展览" sheet, same F/G update
$ws4.Range("F22").Value = 307
$ws4.Range("G22").Value = 218
